$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells: "_old" -> "_FV2404" and "_new" -> "_FV2410" ---
# Columns A..J (1-10) carry the "_old" suffixed headers, column K (11) is "diff",
# and columns L..U (12-21) carry the "_new" suffixed headers.
$leftHeaders = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

$rightHeaders = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $leftHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $leftHeaders[$i]
}

for ($i = 0; $i -lt $rightHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $rightHeaders[$i]
}

# --- 2. Turn the used range into an Excel table ("Table1") ---
$range = $ws.Range("A1:U68")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
# Match the source workbook, which has no named table style applied.
$tbl.TableStyle = ""

# --- 3. Freeze the header row (split below row 1, pane starting at A2) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
